$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# "Running suites c and d": flip the Runmode column (C2:C6) from "N" to "Y"
# so those test-case rows are included in the run.
$ws.Range("C2:C6").Value = "Y"

# Leave the active selection on C6.
$ws.Range("C6").Select()
